# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the "K" column (column G) with recalculated strikeout values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new K value (column G), only rows whose recomputed value changed
$kValues = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 1
    8  = 1
    9  = 2
    10 = 0
    11 = 3
    12 = 0
    13 = 0
    14 = 1
    15 = 3
    16 = 4
    17 = 5
    18 = 2
    19 = 2
    20 = 3
    21 = 1
    23 = 1
    25 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
